# Weekly update: insert two new daily price records for "Pimiento" (Zafiro
# rojo / Zafiro verde) at the top of the data block (rows 292-293), pushing
# the existing historical rows down by two (old row 292 -> new row 294, and
# so on through the end of the sheet, old row 356 -> new row 358).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 292..356 down to 294..358, leaving two fresh blank rows at 292:293.
$ws.Rows("292:293").Insert()

# New record 1: Zafiro rojo
$ws.Range("A292").Value = 7
$ws.Range("B292").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C292").Value = "Ñuble"
$ws.Range("D292").Value = 44889
$ws.Range("E292").Value = 16
$ws.Range("F292").Value = 100112002
$ws.Range("G292").Value = "Pimiento"
$ws.Range("H292").Value = "Zafiro rojo"
$ws.Range("I292").Value = "Primera"
$ws.Range("J292").Value = 120
$ws.Range("K292").Value = 17000
$ws.Range("L292").Value = 18000
$ws.Range("M292").Value = 17500
$ws.Range("N292").Value = "`$/caja 15 kilos"
$ws.Range("O292").Value = "Región de Arica y Parinacota"
$ws.Range("P292").Value = 1167
$ws.Range("Q292").Value = 15
$ws.Range("R292").Value = "Hortaliza"

# New record 2: Zafiro verde
$ws.Range("A293").Value = 7
$ws.Range("B293").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C293").Value = "Ñuble"
$ws.Range("D293").Value = 44889
$ws.Range("E293").Value = 16
$ws.Range("F293").Value = 100112002
$ws.Range("G293").Value = "Pimiento"
$ws.Range("H293").Value = "Zafiro verde"
$ws.Range("I293").Value = "Primera"
$ws.Range("J293").Value = 120
$ws.Range("K293").Value = 15000
$ws.Range("L293").Value = 16000
$ws.Range("M293").Value = 15500
$ws.Range("N293").Value = "`$/caja 15 kilos"
$ws.Range("O293").Value = "Región de Arica y Parinacota"
$ws.Range("P293").Value = 1033
$ws.Range("Q293").Value = 15
$ws.Range("R293").Value = "Hortaliza"
